$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 95, pushing the existing row 95 (and all
# rows below it) down by one. This preserves the formatting of the row
# that used to occupy row 95 (it ends up at row 96 afterwards).
$ws.Rows("95:95").Insert()

# Populate the newly inserted row 95 with the new weekly price record.
# Columns A, B, C, E, F, G, H, I, J, K, L, Q, R, T keep the same values
# as the neighboring Membrillo / Champion / Primera rows.
$ws.Range("A95").Value = 10
$ws.Range("B95").Value = "Vega Modelo de Temuco"
$ws.Range("C95").Value = "La Araucanía"
$ws.Range("D95").Value = 44729
$ws.Range("D95").NumberFormat = $ws.Range("D96").NumberFormat
$ws.Range("E95").Value = 9
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100104
$ws.Range("H95").Value = "Frutos de pepita"
$ws.Range("I95").Value = 100104003
$ws.Range("J95").Value = "Membrillo"
$ws.Range("K95").Value = "Champion"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 140
$ws.Range("N95").Value = 10000
$ws.Range("O95").Value = 10000
$ws.Range("P95").Value = 10000
$ws.Range("Q95").Value = "$/bandeja 18 kilos granel"
$ws.Range("R95").Value = "Región de O'Higgins"
$ws.Range("S95").Value = 556
$ws.Range("T95").Value = 18
